$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 685.9545000000001
$ws.Range("I33").Value = 493.29413
$ws.Range("J33").Value = 1341
$ws.Range("K33").Value = 493.29413
$ws.Range("L33").Value = 1341
$ws.Range("M33").Value = -264.29413
$ws.Range("N33").Value = -1799

$ws.Range("H107").Value = 1598.9
$ws.Range("I107").Value = 939.8823
$ws.Range("J107").Value = 5333.3335
$ws.Range("K107").Value = 939.8823
$ws.Range("L107").Value = 5333.3335
$ws.Range("M107").Value = 980.1177
$ws.Range("N107").Value = -9173.333500000001

$ws.Range("H113").Value = 12003
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 12003
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 12003
$ws.Range("N113").Value = -18511
$ws.Range("M113").ClearContents()

$ws.Range("H138").Value = 3548.9207
$ws.Range("I138").Value = 1732.0513
$ws.Range("J138").Value = 6501.3335
$ws.Range("K138").Value = 5196.1539
$ws.Range("L138").Value = 19504.0005
$ws.Range("M138").Value = -56.15390000000025
$ws.Range("N138").Value = -29784.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10419130
$ws.Range("I2").Value = 35715428
$ws.Range("J2").Value = 3006.647
$ws.Range("K2").Value = 35715428
$ws.Range("L2").Value = 3006.647
$ws.Range("M2").Value = -35715315

$ws.Range("H32").Value = 4942.7847
$ws.Range("I32").Value = 3672.4727
$ws.Range("J32").Value = 11929.5
$ws.Range("K32").Value = 3672.4727
$ws.Range("L32").Value = 11929.5
$ws.Range("M32").Value = -3385.4727
$ws.Range("N32").Value = -12503.5

$ws.Range("H45").Value = 1832.5
$ws.Range("I45").Value = 1062.3334
$ws.Range("J45").Value = 6453.5
$ws.Range("K45").Value = 1062.3334
$ws.Range("L45").Value = 6453.5
$ws.Range("M45").Value = -685.3334
$ws.Range("N45").Value = -7207.5

$ws.Range("H61").Value = 2480.375
$ws.Range("I61").Value = 1104.7273
$ws.Range("J61").Value = 3644.3845
$ws.Range("K61").Value = 1104.7273
$ws.Range("L61").Value = 3644.3845
$ws.Range("M61").Value = -892.7273
$ws.Range("N61").Value = -4068.3845

$ws.Range("H116").Value = 10419130
$ws.Range("I116").Value = 35715428
$ws.Range("J116").Value = 3006.647
$ws.Range("K116").Value = 35715428
$ws.Range("L116").Value = 3006.647
$ws.Range("M116").Value = -35713134

$ws.Range("H122").Value = 3190.4119
$ws.Range("I122").Value = 2191.4443
$ws.Range("J122").Value = 4314.25
$ws.Range("K122").Value = 6574.3329
$ws.Range("L122").Value = 12942.75
$ws.Range("M122").Value = -4124.3329
$ws.Range("N122").Value = -17842.75

$ws.Range("H132").Value = 13160544
$ws.Range("I132").Value = 17546352
$ws.Range("J132").Value = 3120.9473
$ws.Range("K132").Value = 52639056
$ws.Range("L132").Value = 9362.841899999999
$ws.Range("M132").Value = -52636526

$ws.Range("H136").Value = 2480.375
$ws.Range("I136").Value = 1104.7273
$ws.Range("J136").Value = 3644.3845
$ws.Range("K136").Value = 3314.1819
$ws.Range("L136").Value = 10933.1535
$ws.Range("M136").Value = -764.1819
$ws.Range("N136").Value = -16033.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10419130
$ws.Range("I3").Value = 35715428
$ws.Range("J3").Value = 3006.647
$ws.Range("K3").Value = 35715428
$ws.Range("L3").Value = 3006.647
$ws.Range("M3").Value = -35715314

$ws.Range("H107").Value = 2339.3333
$ws.Range("I107").Value = 1346.2
$ws.Range("J107").Value = 3580.75
$ws.Range("K107").Value = 1346.2
$ws.Range("L107").Value = 3580.75
$ws.Range("M107").Value = 573.8
$ws.Range("N107").Value = -7420.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16669693
$ws.Range("I58").Value = 2144.3
$ws.Range("J58").Value = 50004790
$ws.Range("K58").Value = 2144.3
$ws.Range("L58").Value = 50004790
$ws.Range("M58").Value = -1941.3
$ws.Range("N58").Value = -50005196

$ws.Range("H122").Value = 4058.4167
$ws.Range("I122").Value = 3400.2856
$ws.Range("J122").Value = 4979.8
$ws.Range("K122").Value = 10200.8568
$ws.Range("L122").Value = 14939.4
$ws.Range("M122").Value = -7750.856800000001
$ws.Range("N122").Value = -19839.4

$ws.Range("H132").Value = 4220.28
$ws.Range("I132").Value = 2928.0908
$ws.Range("J132").Value = 5235.5713
$ws.Range("K132").Value = 8784.2724
$ws.Range("L132").Value = 15706.7139
$ws.Range("M132").Value = -6254.2724
$ws.Range("N132").Value = -20766.7139

$ws.Range("H136").Value = 16669693
$ws.Range("I136").Value = 2144.3
$ws.Range("J136").Value = 50004790
$ws.Range("K136").Value = 6432.900000000001
$ws.Range("L136").Value = 150014370
$ws.Range("M136").Value = -3882.900000000001
$ws.Range("N136").Value = -150019470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 10672.454
$ws.Range("I34").Value = 200
$ws.Range("J34").Value = 12999.667
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 38999.001
$ws.Range("M34").Value = -516
$ws.Range("N34").Value = -39167.001

$ws.Range("H39").Value = 2975
$ws.Range("I39").Value = 750
$ws.Range("J39").Value = 3716.6667
$ws.Range("K39").Value = 2250
$ws.Range("L39").Value = 11150.0001
$ws.Range("M39").Value = -1956
$ws.Range("N39").Value = -11738.0001

$ws.Range("H55").Value = 2974.8
$ws.Range("I55").Value = 475
$ws.Range("J55").Value = 3599.75
$ws.Range("K55").Value = 1425
$ws.Range("L55").Value = 10799.25
$ws.Range("M55").Value = -1248
$ws.Range("N55").Value = -11153.25

$ws.Range("H57").Value = 500
$ws.Range("I57").Value = 500
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 1500
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -941
$ws.Range("N57").ClearContents()

$ws.Range("H96").Value = 3800
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 3800
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 11400
$ws.Range("N96").Value = -15518

$ws.Range("H98").Value = 100
$ws.Range("I98").Value = 100
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 300
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 1198
$ws.Range("N98").ClearContents()

$ws.Range("H113").Value = 703.2
$ws.Range("I113").Value = 478.44446
$ws.Range("J113").Value = 829.625
$ws.Range("K113").Value = 1435.33338
$ws.Range("L113").Value = 2488.875
$ws.Range("M113").Value = 734.66662
$ws.Range("N113").Value = -6828.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3794.9722
$ws.Range("I132").Value = 2817.5417
$ws.Range("J132").Value = 5749.8335
$ws.Range("K132").Value = 8452.625100000001
$ws.Range("L132").Value = 17249.5005
$ws.Range("M132").Value = -5922.625100000001
$ws.Range("N132").Value = -22309.5005

$ws.Range("H135").Value = 22097.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 22097.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 22097.5
$ws.Range("N135").Value = -32237.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 142859600
$ws.Range("I22").Value = 333333820
$ws.Range("J22").Value = 3924.75
$ws.Range("K22").Value = 333333820
$ws.Range("L22").Value = 3924.75
$ws.Range("M22").Value = -333333525

$ws.Range("H27").Value = 142859600
$ws.Range("I27").Value = 333333820
$ws.Range("J27").Value = 3924.75
$ws.Range("K27").Value = 333333820
$ws.Range("L27").Value = 3924.75
$ws.Range("M27").Value = -333333713

$ws.Range("H132").Value = 4421.5713
$ws.Range("I132").Value = 3300.8
$ws.Range("J132").Value = 5044.222
$ws.Range("K132").Value = 9902.400000000001
$ws.Range("L132").Value = 15132.666
$ws.Range("M132").Value = -7372.400000000001
$ws.Range("N132").Value = -20192.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 261940.89
$ws.Range("I132").Value = 402207.88
$ws.Range("J132").Value = 11464.143
$ws.Range("K132").Value = 1206623.64
$ws.Range("L132").Value = 34392.429
$ws.Range("M132").Value = -1204093.64
$ws.Range("N132").Value = -39452.429
